$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the empty "_GoBack" bookmark in the "Events Loose Ideas"
#    heading paragraph, without disturbing the existing run split.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Italicise the "The texts below are loose ideas..." paragraph
#    (both the run and the paragraph mark).
# ------------------------------------------------------------------
$introPara = $d.Paragraphs(4)
$introRange = $introPara.Range
$introRange.Font.Italic = $true
$introRange.Font.ItalicBi = $true

# ------------------------------------------------------------------
# 3. Fix "Ofcourse" -> "Of" + " " + "course" (three separate runs)
#    and wrap the "lne" typo with spell-check proofErr markers, by
#    replacing the paragraph's run content with the equivalent OOXML.
# ------------------------------------------------------------------
$target = "Ofcourse you can easily disambiguate: either set the rule that comment has an inward lne, while events have an outward line. Notation however should go regardless of the placement of a name, so if you ignore the names, the notations do indeed conflict. Then there is the second option to disambiguation: always show the public of friend access mark and in case of private either show an explicit private access mark, or for an explicitly private notation, keep the event connector line inside the container (the lower-left notation)."

$found = $d.Content
$found.Find.Text = $target
$found.Find.Execute() | Out-Null
$matchRange = $d.Range($found.Start, $found.End)

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Of</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">course you can easily disambiguate: either set the rule that comment has an inward </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lne</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, while events have an outward line. Notation however should go regardless of the placement of a name, so if you ignore the names, the notations do indeed conflict. Then there is the second option to disambiguation: always show the public of friend access mark and in case of private either show an explicit private access mark, or for an explicitly private notation, keep the event connector line inside the container (the lower-left notation).</w:t></w:r></w:p>'

$matchRange.InsertXML($newXml)
